$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells are treated as text, matching the source workbook
# which stores these values as inline strings (not numbers).
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '303.29'
$ws.Range("E2").Value = '5.61%'
$ws.Range("D3").Value = '31.92'
$ws.Range("E3").Value = '9.43%'
$ws.Range("D4").Value = '5.267'
$ws.Range("E4").Value = '1.33%'
$ws.Range("D5").Value = '0.07447'
$ws.Range("E5").Value = '6.32%'
$ws.Range("D6").Value = '7.855'
$ws.Range("E6").Value = '5.64%'
$ws.Range("D7").Value = '3.800'
$ws.Range("E7").Value = '7.00%'
$ws.Range("D8").Value = '1.518'
$ws.Range("E8").Value = '7.70%'
$ws.Range("D9").Value = '0.9191'
$ws.Range("E9").Value = '1.86%'
$ws.Range("D10").Value = '0.01754'
$ws.Range("E10").Value = '2,602.30%'
$ws.Range("D11").Value = '0.1687'
$ws.Range("E11").Value = '4.99%'
$ws.Range("D12").Value = '0.08031'
$ws.Range("E12").Value = '7.86%'
$ws.Range("D13").Value = '0.07941'
$ws.Range("E13").Value = '2.81%'
$ws.Range("D14").Value = '0.03047'
$ws.Range("E14").Value = '4.12%'
$ws.Range("D15").Value = '0.09893'
$ws.Range("E15").Value = '9.69%'
$ws.Range("D16").Value = '0.001519'
$ws.Range("E16").Value = '-4.13%'
$ws.Range("D17").Value = '0.04619'
$ws.Range("E17").Value = '2.09%'
$ws.Range("D18").Value = '0.006150'
$ws.Range("E18").Value = '0.80%'
$ws.Range("D19").Value = '3.472'
$ws.Range("E19").Value = '0.12%'
$ws.Range("E20").Value = '0.01%'
$ws.Range("E21").Value = '2.65%'
$ws.Range("D22").Value = '0.1330'
$ws.Range("E22").Value = '-0.18%'
$ws.Range("D23").Value = '4.485'
$ws.Range("E23").Value = '11.97%'
$ws.Range("D24").Value = '0.1625'
$ws.Range("E24").Value = '1.67%'
$ws.Range("D25").Value = '0.001222'
$ws.Range("E25").Value = '1.01%'
$ws.Range("D26").Value = '0.004439'
$ws.Range("E26").Value = '4.62%'
$ws.Range("D27").Value = '0.0001400'
$ws.Range("E27").Value = '19.82%'
$ws.Range("D28").Value = '0.0001749'
$ws.Range("E28").Value = '4.92%'
$ws.Range("D40").Value = '0.04490'
$ws.Range("E40").Value = '3.26%'
$ws.Range("D41").Value = '0.007172'
$ws.Range("E41").Value = '3.12%'
$ws.Range("D42").Value = '0.1348'
$ws.Range("E42").Value = '8.10%'
$ws.Range("D43").Value = '0.002210'
$ws.Range("E43").Value = '6.91%'
$ws.Range("D44").Value = '0.01277'
$ws.Range("E44").Value = '9.61%'
$ws.Range("D45").Value = '0.00006160'
$ws.Range("E45").Value = '5.75%'
$ws.Range("D46").Value = '1.868'
$ws.Range("E46").Value = '-3.16%'
$ws.Range("E47").Value = '-0.25%'
